# "Fixes sort part 2 - capitalization"
# Column A holds city names that (aside from the header) were stored in
# mixed/title case while every other column (B:BG) already used the
# all-caps city codes. Re-key column A to match that same all-caps
# convention, and rename the header from "No" to "City".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'City'
$ws.Range("A2").Value = 'ΛΕΥΚΩΣΙΑ'
$ws.Range("A3").Value = 'ΚΑΙΜΑΚΛΙ'
$ws.Range("A4").Value = 'ΑΓΛΑΝΤΖΙΑ'
$ws.Range("A5").Value = 'ΕΓΚΩΜΗ'
$ws.Range("A6").Value = 'ΑΓ. ΤΡΙΜΙΘΙΑΣ'
$ws.Range("A7").Value = 'ΑΚΑΚΙ'
$ws.Range("A8").Value = 'ΑΛΑΜΠΡΑ'
$ws.Range("A9").Value = 'ΝΗΣΟΥ'
$ws.Range("A10").Value = 'ΑΣΤΡΟΜΕΡΙΤΗΣ'
$ws.Range("A11").Value = 'ΓΕΡΙ'
$ws.Range("A12").Value = 'ΔΑΛΙ'
$ws.Range("A13").Value = 'ΔΕΥΤΕΡΑ'
$ws.Range("A14").Value = 'ΕΡΓΑΤΕΣ'
$ws.Range("A15").Value = 'ΕΥΡΥΧΟΥ'
$ws.Range("A16").Value = 'ΚΟΚΚΙΝΟΤΡΙΜΙΘΙΑ'
$ws.Range("A17").Value = 'ΛΑΚΑΤΑΜΙΑ'
$ws.Range("A18").Value = 'ΛΑΤΣΙΑ'
$ws.Range("A19").Value = 'ΛΥΜΠΙΑ'
$ws.Range("A20").Value = 'ΜΑΘΙΑΤΗΣ'
$ws.Range("A21").Value = 'ΜΑΜΜΑΡΙ'
$ws.Range("A22").Value = 'ΜΕΝΟΙΚΟ'
$ws.Range("A23").Value = 'ΛΥΘΡΟΔΟΝΤΑΣ'
$ws.Range("A24").Value = 'ΤΣΕΡΙ'
$ws.Range("A25").Value = 'ΨΗΜΟΛΟΦΟΥ'
$ws.Range("A26").Value = 'ΛΕΜΕΣΟΣ'
$ws.Range("A27").Value = 'ΚΟΛΟΣΣΙ'
$ws.Range("A28").Value = 'ΑΛΑΣΣΑ'
$ws.Range("A29").Value = 'ΑΓ. ΙΩΑΝΝΗΣ'
$ws.Range("A30").Value = 'ΥΨΩΝΑΣ'
$ws.Range("A31").Value = 'ΛΑΡΝΑΚΑ'
$ws.Range("A32").Value = 'ΑΘΗΑΙΝΟΥ'
$ws.Range("A33").Value = 'ΚΙΤΙ'
$ws.Range("A34").Value = 'ΚΟΡΝΟΣ'
$ws.Range("A35").Value = 'ΑΛΕΘΡΙΚΟ'
$ws.Range("A36").Value = 'ΚΟΦΙΝΟΥ'
$ws.Range("A37").Value = 'ΨΕΥΔΑΣ'
$ws.Range("A38").Value = 'ΛΕΙΒΑΔΕΙΑ'
$ws.Range("A39").Value = 'ΞΥΛΟΦΑΓΟΥ'
$ws.Range("A40").Value = 'ΖΥΓΙ'
$ws.Range("A41").Value = 'ΠΥΛΑ'
$ws.Range("A42").Value = 'ΠΑΡΑΛΙΜΝΙ'
$ws.Range("A43").Value = 'ΑΥΓΟΡΟΥ'
$ws.Range("A44").Value = 'ΛΙΟΠΕΤΡΙ'
$ws.Range("A45").Value = 'ΣΩΤΗΡΑ'
$ws.Range("A46").Value = 'ΦΡΕΝΑΡΟΣ'
$ws.Range("A47").Value = 'ΠΑΦΟΣ'
$ws.Range("A48").Value = 'ΑΝΑΡΙΤΑ'
$ws.Range("A49").Value = 'ΕΜΠΑ'
$ws.Range("A50").Value = 'ΠΕΓΕΙΑ'
$ws.Range("A51").Value = 'ΑΓ. ΜΑΡΙΝΑ ΧΡΥΣ.'
$ws.Range("A52").Value = 'ΠΟΛΗ ΧΡΥΣΟΧΟΥΣ'
$ws.Range("A53").Value = 'ΤΣΑΔΑ'

# Restore the default selection (top-left cell) instead of the stale
# B1:BG1 selection left over from the previous save.
$ws.Range("A1").Select()
